# "add coffee pest and disease mgt"
# The source data this report was generated from dropped one row
# (Kitagwenda / Nganiko / Kagorogoro_A / Pruning saw / 20). Remove the
# corresponding row from the worksheet so everything below shifts up,
# matching the refreshed export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("22").Delete()
